$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 5 (DVC_SAMPLE_DATA), pushing it down to row 10
$ws.Range("A5:E9").EntireRow.Insert()

# Row 5: DVC_2PM_DCRESISTANCE_1
$ws.Range("A5").Value = "DVC_2PM_DCRESISTANCE_1"
$ws.Range("B5").Value = "23-0: Test_Current_Value"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "RW"

# Row 6: DVC_2PM_CURRVOLT_1
$ws.Range("A6").Value = "DVC_2PM_CURRVOLT_1"
$ws.Range("B6").Value = "1-0: Sweep_Param"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "RW"

# Row 7: DVC_2PM_CURRVOLT_2
$ws.Range("A7").Value = "DVC_2PM_CURRVOLT_2"
$ws.Range("B7").Value = "23-0: Starting_Param"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = "RW"

# Row 8: DVC_2PM_CURRVOLT_3
$ws.Range("A8").Value = "DVC_2PM_CURRVOLT_3"
$ws.Range("B8").Value = "23-0: Ending_Param"
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = "RW"

# Row 9: DVC_2PM_CURRVOLT_4
$ws.Range("A9").Value = "DVC_2PM_CURRVOLT_4"
$ws.Range("B9").Value = "23-0: Increment_Param"
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = "RW"

# Apply row heights matching the target
$ws.Range("A5:E9").RowHeight = 13.2

# Update selection to match target view
$ws.Range("G6:G8").Select()
